$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "43.699.99"
Set-TextValue "D3" "2.319.10"
Set-TextValue "E3" "  +4.76%  "
Set-TextValue "E4" "  -0.05%  "
Set-TextValue "D5" "269.22"
Set-TextValue "E5" "  -0.18%  "
Set-TextValue "D6" "92.02"
Set-TextValue "E6" "  +7.44%  "
Set-TextValue "D7" "0.631"
Set-TextValue "E7" "  +1.76%  "
Set-TextValue "E8" "  +0.03%  "
Set-TextValue "D9" "0.621"
Set-TextValue "E9" "  +3.09%  "
Set-TextValue "D10" "44.80"
Set-TextValue "E10" "  -2.00%  "
Set-TextValue "D11" "0.0936"
Set-TextValue "E11" "  +1.78%  "
Set-TextValue "D12" "8.00"
Set-TextValue "E12" "  +6.93%  "
Set-TextValue "E13" "  +0.23%  "
Set-TextValue "D14" "2.666.74"
Set-TextValue "E14" "  +4.79%  "
Set-TextValue "D15" "15.29"
Set-TextValue "E15" "  +5.00%  "
Set-TextValue "D16" "0.855"
Set-TextValue "E16" "  +9.20%  "
Set-TextValue "D17" "2.321.94"
Set-TextValue "E17" "  +4.77%  "
Set-TextValue "D18" "43.678.76"
Set-TextValue "E18" "  +0.06%  "
Set-TextValue "E19" "  +2.99%  "
Set-TextValue "D20" "6.32"
Set-TextValue "E20" "  +5.72%  "
Set-TextValue "D21" "71.29"
Set-TextValue "E21" "  +2.09%  "
Set-TextValue "D22" "241.42"
Set-TextValue "E22" "  +4.28%  "
Set-TextValue "D23" "2.28"
Set-TextValue "E23" "  -3.53%  "
Set-TextValue "D24" "9.67"
Set-TextValue "E24" "  +9.43%  "
Set-TextValue "D25" "0.999"
Set-TextValue "E25" "  -0.06%  "
Set-TextValue "E26" "  -8.22%  "
Set-TextValue "E27" "  +4.44%  "
Set-TextValue "E28" "  +5.08%  "
Set-TextValue "D29" "3.37"
Set-TextValue "E29" "  -4.48%  "
Set-TextValue "D30" "38.84"
Set-TextValue "E30" "  -0.32%  "
Set-TextValue "D31" "22.52"
Set-TextValue "E31" "  +9.80%  "
Set-TextValue "D32" "172.76"
Set-TextValue "E32" "  -1.47%  "
Set-TextValue "D33" "0.0891"
Set-TextValue "E33" "  +0.16%  "
Set-TextValue "D34" "5.54"
Set-TextValue "E34" "  +2.98%  "
Set-TextValue "E35" "  +1.69%  "
Set-TextValue "E36" "  +0.78%  "
Set-TextValue "D37" "4.50"
Set-TextValue "E37" "  +3.47%  "
Set-TextValue "E38" "  -1.79%  "
Set-TextValue "D39" "3.36"
Set-TextValue "E39" "  +3.15%  "
Set-TextValue "D40" "0.236"
Set-TextValue "E40" "  +16.13%  "
Set-TextValue "D41" "2.31"
Set-TextValue "E41" "  +10.64%  "
Set-TextValue "D42" "12.24"
Set-TextValue "E42" "  -0.08%  "
Set-TextValue "E43" "  +19.26%  "
Set-TextValue "E44" "  +0.76%  "
Set-TextValue "D45" "61.24"
Set-TextValue "E45" "  -5.73%  "
Set-TextValue "D46" "8.92"
Set-TextValue "E46" "  +7.16%  "
Set-TextValue "D47" "0.102"
Set-TextValue "E47" "  +2.99%  "
Set-TextValue "D48" "100.53"
Set-TextValue "E48" "  +0.33%  "
Set-TextValue "E49" "  -0.46%  "
Set-TextValue "D50" "2.543.72"
Set-TextValue "E50" "  +4.73%  "
Set-TextValue "D51" "0.431"
Set-TextValue "E51" "  -1.65%  "
